$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (squad)
$ws.Range("A2").Value = "Mayorista"
$ws.Range("A3").Value = "Mayorista"
$ws.Range("A4").Value = "Mayorista"
$ws.Range("A5").Value = "Negocio"
$ws.Range("A6").Value = "Negocio"

# Column B (puesto)
$ws.Range("B2").Value = "Ejecutivo de Negocio"
$ws.Range("B3").Value = "Analista de Negocio"
$ws.Range("B4").Value = "Jefe de Riesgo"
$ws.Range("B5").Value = "Funcionario de Negocio"
$ws.Range("B6").Value = "Po - Banca Negocio"

# Column C (nombreUsuario)
$ws.Range("C2").Value = "Alex Mejia"
$ws.Range("C3").Value = "Susana Flores"
$ws.Range("C4").Value = "Juan Ruiz"
$ws.Range("C5").Value = "Pedro Montex"
$ws.Range("C6").Value = "Jorge Olivares"

# Column D (manager) - only rows 2 and 4 keep a value now
$ws.Range("D2").Value = "Juan Ruiz"
$ws.Range("D3").ClearContents()
$ws.Range("D4").Value = "Manager de prueba"
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()

# Column E (delegado) - only row 6 keeps a value now
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").Value = "Iveth Mattos"

# Widen column B to fit the new, longer "puesto" labels
$ws.Columns("B").AutoFit()

# Move the active selection to F2 (matches the saved sheet view)
$ws.Range("F2").Select()
